$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row of data for 2020-08-03 (raw and clean SSA data)
$newRow = 65

# Force column A to be entered as text so "2020-08-03" isn't auto-converted
# to a date serial number, then clear the temporary formatting so the cell
# ends up with the same (default) style as the rest of the date column.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2020-08-03"
$ws.Cells.Item($newRow, 1).ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 443813
$ws.Cells.Item($newRow, 3).Value = 488207
$ws.Cells.Item($newRow, 4).Value = 79030
$ws.Cells.Item($newRow, 5).Value = 48012
$ws.Cells.Item($newRow, 6).Value = 26.88
